$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 6.0.0 -> 6.1.0
$meta.Range("B3").Value = "6.1.0"

# Date: 2022-01-21T20:46:54+00:00 -> 2022-05-31T20:10:14+00:00
$meta.Range("B8").Value = "2022-05-31T20:10:14+00:00"

# --- Concepts sheet updates ---
$concepts = $wb.Worksheets.Item("Concepts")

# Row 40: np-std-outreach-pending -> np-std-care-gap-pending
$concepts.Range("B40").Value = "np-std-care-gap-pending"
$concepts.Range("C40").Value = "Not Processed - Standard care gap pending"
$concepts.Range("D40").Value = "CommunicationRequest not processed because conversation was not activated after multiple attempts. Standard care gap communication sent."
